# Update "想去人数" (want-to-go count) values on two sheets to reflect
# the latest scrape, as generated by the gh-pages build at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 647
$wsExhibit.Range("F4").Value = 1476
$wsExhibit.Range("F5").Value = 689

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 647
$wsAll.Range("F4").Value = 1476
$wsAll.Range("F6").Value = 689
